# Revert project to commit f4eed51
# Restore the original (pre-edit) ordering / values of the link-check rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> google.com, checked, success
$ws.Range("A2").Value = "https://www.google.com"
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = "Checked"
$ws.Range("D2").Value = "Success"

# Row 3 -> sebi.gov.in (status unchanged: 200 / Checked / Success)
$ws.Range("A3").Value = "https://www.sebi.gov.in"

# Row 4 -> surveymonkey.com (unchanged)

# Row 5 -> owasp.org (status unchanged: 200 / Checked / Success)
$ws.Range("A5").Value = "https://www.owasp.org"

# Row 6 -> axisbank.com (status unchanged: 200 / Checked / Success)
$ws.Range("A6").Value = "https://www.axisbank.com"

# Row 7 -> sahilendworldfibvweuidbuk.org, not checked, connection error
$ws.Range("A7").Value = "https://www.sahilendworldfibvweuidbuk.org"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "Not Checked"
$ws.Range("D7").Value = "HTTPSConnectionPool(host='www.sahilendworldfibvweuidbuk.org', port=443): Max retries exceeded with url: / (Caused by NameResolutionError(`"<urllib3.connection.HTTPSConnection object at 0x000001CC7FEB4550>: Failed to resolve 'www.sahilendworldfibvweuidbuk.org' ([Errno 11001] getaddrinfo failed)`"))"

# Row 8 -> rbi.org.in (status unchanged: 200 / Checked / Success)
$ws.Range("A8").Value = "https://www.rbi.org.in"
